# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (F) and "最低票价" (G) numeric counters across sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / rId1)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 432
$wsExpo.Range("F3").Value = 2966
$wsExpo.Range("F5").Value = 71

# Sheet "演出" (sheet2 / rId2)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("G6").Value = 72

# Sheet "全部类型" (sheet4 / rId4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 432
$wsAll.Range("F7").Value = 2966
$wsAll.Range("G9").Value = 72
$wsAll.Range("F10").Value = 71
